$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ntn1"
$ws.Range("C2").Value = "Adora2b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.9305633333333333
$ws.Range("H2").Value = 2.79169
$ws.Range("I2").Value = 0.01768777137856805
$ws.Range("J2").Value = 0.01768777137856806
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.237097333333333
$ws.Range("N2").Value = 3.711292
$ws.Range("O2").Value = 0.1496639399539542
$ws.Range("P2").Value = 0.1496639399539542
$ws.Range("Q2").Value = 1.151197418164445
$ws.Range("R2").Value = 10.36077676348
$ws.Range("S2").Value = 0.00264722155352128
$ws.Range("T2").Value = 0.00264722155352128

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ntn1"
$ws.Range("C3").Value = "Adora2b"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.9305633333333333
$ws.Range("H3").Value = 2.79169
$ws.Range("I3").Value = 0.01768777137856805
$ws.Range("J3").Value = 0.01768777137856806
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.467117666666667
$ws.Range("N3").Value = 4.401353
$ws.Range("O3").Value = 0.1774917821362901
$ws.Range("P3").Value = 0.1774917821362901
$ws.Range("Q3").Value = 1.365245906285556
$ws.Range("R3").Value = 12.28721315657
$ws.Range("S3").Value = 0.003139434064001308
$ws.Range("T3").Value = 0.003139434064001309

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntn1"
$ws.Range("C4").Value = "Adora2b"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9305633333333333
$ws.Range("H4").Value = 2.79169
$ws.Range("I4").Value = 0.01768777137856805
$ws.Range("J4").Value = 0.01768777137856806
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.561619333333333
$ws.Range("N4").Value = 16.684858
$ws.Range("O4").Value = 0.6728442779097555
$ws.Range("P4").Value = 0.6728442779097557
$ws.Range("Q4").Value = 5.175439025557777
$ws.Range("R4").Value = 46.57895123002
$ws.Range("S4").Value = 0.01190111576104546
$ws.Range("T4").Value = 0.01190111576104547

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Ntn1"
$ws.Range("C5").Value = "Adora2b"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 44.154177
$ws.Range("H5").Value = 132.462531
$ws.Range("I5").Value = 0.8392647337471152
$ws.Range("J5").Value = 0.8392647337471153
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.237097333333333
$ws.Range("N5").Value = 3.711292
$ws.Range("O5").Value = 0.1496639399539542
$ws.Range("P5").Value = 0.1496639399539542
$ws.Range("Q5").Value = 54.62301462222801
$ws.Range("R5").Value = 491.6071316000521
$ws.Range("S5").Value = 0.1256076667169996
$ws.Range("T5").Value = 0.1256076667169997

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Ntn1"
$ws.Range("C6").Value = "Adora2b"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 44.154177
$ws.Range("H6").Value = 132.462531
$ws.Range("I6").Value = 0.8392647337471152
$ws.Range("J6").Value = 0.8392647337471153
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.467117666666667
$ws.Range("N6").Value = 4.401353
$ws.Range("O6").Value = 0.1774917821362901
$ws.Range("P6").Value = 0.1774917821362901
$ws.Range("Q6").Value = 64.77937313382701
$ws.Range("R6").Value = 583.014358204443
$ws.Range("S6").Value = 0.1489625932769145
$ws.Range("T6").Value = 0.1489625932769145

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Ntn1"
$ws.Range("C7").Value = "Adora2b"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 44.154177
$ws.Range("H7").Value = 132.462531
$ws.Range("I7").Value = 0.8392647337471152
$ws.Range("J7").Value = 0.8392647337471153
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.561619333333333
$ws.Range("N7").Value = 16.684858
$ws.Range("O7").Value = 0.6728442779097555
$ws.Range("P7").Value = 0.6728442779097557
$ws.Range("Q7").Value = 245.568724450622
$ws.Range("R7").Value = 2210.118520055598
$ws.Range("S7").Value = 0.5646944737532009
$ws.Range("T7").Value = 0.5646944737532011

# Row 8
$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Ntn1"
$ws.Range("C8").Value = "Adora2b"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.525807
$ws.Range("H8").Value = 22.577421
$ws.Range("I8").Value = 0.1430474948743168
$ws.Range("J8").Value = 0.1430474948743168
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.237097333333333
$ws.Range("N8").Value = 3.711292
$ws.Range("O8").Value = 0.1496639399539542
$ws.Range("P8").Value = 0.1496639399539542
$ws.Range("Q8").Value = 9.310155770881336
$ws.Range("R8").Value = 83.79140193793201
$ws.Range("S8").Value = 0.02140905168343333
$ws.Range("T8").Value = 0.02140905168343333

# Row 9
$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Ntn1"
$ws.Range("C9").Value = "Adora2b"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.525807
$ws.Range("H9").Value = 22.577421
$ws.Range("I9").Value = 0.1430474948743168
$ws.Range("J9").Value = 0.1430474948743168
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.467117666666667
$ws.Range("N9").Value = 4.401353
$ws.Range("O9").Value = 0.1774917821362901
$ws.Range("P9").Value = 0.1774917821362901
$ws.Range("Q9").Value = 11.04124440562367
$ws.Range("R9").Value = 99.37119965061301
$ws.Range("S9").Value = 0.02538975479537431
$ws.Range("T9").Value = 0.02538975479537431

# Row 10
$ws.Range("A10").Value = "ECs"
$ws.Range("B10").Value = "Ntn1"
$ws.Range("C10").Value = "Adora2b"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.525807
$ws.Range("H10").Value = 22.577421
$ws.Range("I10").Value = 0.1430474948743168
$ws.Range("J10").Value = 0.1430474948743168
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.561619333333333
$ws.Range("N10").Value = 16.684858
$ws.Range("O10").Value = 0.6728442779097555
$ws.Range("P10").Value = 0.6728442779097557
$ws.Range("Q10").Value = 41.85567371013533
$ws.Range("R10").Value = 376.701063391218
$ws.Range("S10").Value = 0.09624868839550914
$ws.Range("T10").Value = 0.09624868839550915

Write-Output "Edit complete"
